$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 272-273, shifting the existing rows 272-319 down to 274-321
$ws.Range("A272:A273").EntireRow.Insert()

# Row 272: Early Burlat, Primera
$ws.Range("A272").Value = 5
$ws.Range("B272").Value = "Macroferia Regional de Talca"
$ws.Range("C272").Value = "Maule"
$ws.Range("D272").Value = 45244
$ws.Range("E272").Value = 7
$ws.Range("F272").Value = "Fruta"
$ws.Range("G272").Value = 100103
$ws.Range("H272").Value = "Frutos de hueso (carozo)"
$ws.Range("I272").Value = 100103001
$ws.Range("J272").Value = "Cereza"
$ws.Range("K272").Value = "Early Burlat"
$ws.Range("L272").Value = "Primera"
$ws.Range("M272").Value = 100
$ws.Range("N272").Value = 30000
$ws.Range("O272").Value = 30000
$ws.Range("P272").Value = 30000
$ws.Range("Q272").Value = "`$/bandeja 10 kilos"
$ws.Range("R272").Value = "Provincia de Curicó"
$ws.Range("S272").Value = 3000
$ws.Range("T272").Value = 10

# Row 273: Lapins, Primera
$ws.Range("A273").Value = 5
$ws.Range("B273").Value = "Macroferia Regional de Talca"
$ws.Range("C273").Value = "Maule"
$ws.Range("D273").Value = 45244
$ws.Range("E273").Value = 7
$ws.Range("F273").Value = "Fruta"
$ws.Range("G273").Value = 100103
$ws.Range("H273").Value = "Frutos de hueso (carozo)"
$ws.Range("I273").Value = 100103001
$ws.Range("J273").Value = "Cereza"
$ws.Range("K273").Value = "Lapins"
$ws.Range("L273").Value = "Primera"
$ws.Range("M273").Value = 180
$ws.Range("N273").Value = 30000
$ws.Range("O273").Value = 30000
$ws.Range("P273").Value = 30000
$ws.Range("Q273").Value = "`$/bandeja 10 kilos"
$ws.Range("R273").Value = "Provincia de Curicó"
$ws.Range("S273").Value = 3000
$ws.Range("T273").Value = 10
